$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refresh "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 22:03"

# --- Estados Unidos (row 4): updated totals ---
$ws.Range("B4").Value = 1252470
$ws.Range("C4").Value = 14837
$ws.Range("D4").Value = 205205
$ws.Range("E4").Value = 973473
$ws.Range("G4").Value = 1521
$ws.Range("H4").Value = 73792

# --- Ecuador (row 21): updated totals ---
$ws.Range("E21").Value = 26830
$ws.Range("F21").Value = 156
$ws.Range("G21").Value = 49
$ws.Range("H21").Value = 1618

# --- Suiza (row 22): updated totals ---
$ws.Range("E22").Value = 2855
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 1805

# --- Irlanda (row 28): updated totals ---
$ws.Range("D28").Value = 17110
$ws.Range("E28").Value = 3763

# --- Seychelles / Montserrat: swap order + per-row stats ---
# Seychelles now appears before Montserrat (row 205 becomes Seychelles,
# row 206 becomes Montserrat), each carrying its own updated D/F/H values.
$ws.Range("A205").Value = "Seychelles"
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0

$ws.Range("A206").Value = "Montserrat"
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1
